$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 357.8125
$ws.Range("I5").Value = 215.625
$ws.Range("J5").Value = 500
$ws.Range("K5").Value = 215.625
$ws.Range("L5").Value = 500
$ws.Range("M5").Value = -100.625
$ws.Range("N5").Value = -730
$ws.Range("H128").Value = 41884.285
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 41884.285
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 41884.285
$ws.Range("M128").ClearContents()
$ws.Range("N128").Value = -51844.285
$ws.Range("H141").Value = 33510.938
$ws.Range("I141").Value = 49469.76
$ws.Range("J141").Value = 3044.0908
$ws.Range("K141").Value = 148409.28
$ws.Range("L141").Value = 9132.2724
$ws.Range("M141").Value = -143229.28
$ws.Range("N141").Value = -19492.2724

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 9000
$ws.Range("I3").Value = 9000
$ws.Range("K3").Value = 9000
$ws.Range("M3").Value = -8885
$ws.Range("H15").Value = 3305.2
$ws.Range("J15").Value = 3256.5
$ws.Range("L15").Value = 3256.5
$ws.Range("N15").Value = -3956.5
$ws.Range("H61").Value = 1042.2778
$ws.Range("I61").Value = 986.34283
$ws.Range("K61").Value = 986.34283
$ws.Range("M61").Value = -774.34283
$ws.Range("H136").Value = 1042.2778
$ws.Range("I136").Value = 986.34283
$ws.Range("K136").Value = 2959.02849
$ws.Range("M136").Value = -409.0284900000001

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2991.76
$ws.Range("I134").Value = 1115.7084
$ws.Range("J134").Value = 4723.5
$ws.Range("K134").Value = 3347.1252
$ws.Range("L134").Value = 14170.5
$ws.Range("M134").Value = -812.1251999999999
$ws.Range("N134").Value = -19240.5

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("H7").Value = 285.1875
$ws.Range("I7").Value = 231.08333
$ws.Range("J7").Value = 447.5
$ws.Range("K7").Value = 231.08333
$ws.Range("L7").Value = 447.5
$ws.Range("M7").Value = -118.08333
$ws.Range("N7").Value = -673.5
$ws.Range("H31").Value = 235423.12
$ws.Range("I31").Value = 423290.16
$ws.Range("K31").Value = 423290.16
$ws.Range("M31").Value = -422995.16
$ws.Range("H34").Value = 235423.12
$ws.Range("I34").Value = 423290.16
$ws.Range("K34").Value = 423290.16
$ws.Range("M34").Value = -423088.16
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
$ws.Range("H58").Value = 2586.4707
$ws.Range("I58").Value = 1381.6
$ws.Range("J58").Value = 5933.3335
$ws.Range("K58").Value = 1381.6
$ws.Range("L58").Value = 5933.3335
$ws.Range("M58").Value = -1178.6
$ws.Range("N58").Value = -6339.3335
$ws.Range("H132").Value = 4499.5835
$ws.Range("I132").Value = 3569.2942
$ws.Range("J132").Value = 6758.857
$ws.Range("K132").Value = 10707.8826
$ws.Range("L132").Value = 20276.571
$ws.Range("M132").Value = -8177.882599999999
$ws.Range("N132").Value = -25336.571
$ws.Range("H134").Value = 5607.875
$ws.Range("I134").Value = 5199.522
$ws.Range("K134").Value = 15598.566
$ws.Range("M134").Value = -13063.566
$ws.Range("H136").Value = 2586.4707
$ws.Range("I136").Value = 1381.6
$ws.Range("J136").Value = 5933.3335
$ws.Range("K136").Value = 4144.799999999999
$ws.Range("L136").Value = 17800.0005
$ws.Range("M136").Value = -1594.799999999999
$ws.Range("N136").Value = -22900.0005
$ws.Range("H137").Value = 45508
$ws.Range("J137").Value = 45508
$ws.Range("L137").Value = 45508
$ws.Range("N137").Value = -55708

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3101.3572
$ws.Range("I3").Value = 2428.9
$ws.Range("J3").Value = 4782.5
$ws.Range("K3").Value = 7286.700000000001
$ws.Range("L3").Value = 14347.5
$ws.Range("M3").Value = -7174.700000000001
$ws.Range("N3").Value = -14571.5
$ws.Range("H5").Value = 2572.4211
$ws.Range("I5").Value = 970.5454999999999
$ws.Range("K5").Value = 2911.6365
$ws.Range("M5").Value = -2799.6365
$ws.Range("H12").Value = 136.4
$ws.Range("J12").Value = 191.57143
$ws.Range("L12").Value = 574.71429
$ws.Range("N12").Value = -920.71429
$ws.Range("H113").Value = 491.1143
$ws.Range("I113").Value = 501.10526
$ws.Range("J113").Value = 479.25
$ws.Range("K113").Value = 1503.31578
$ws.Range("L113").Value = 1437.75
$ws.Range("M113").Value = 666.6842200000001
$ws.Range("N113").Value = -5777.75
$ws.Range("H131").Value = 763.61
$ws.Range("I131").Value = 393.5
$ws.Range("J131").Value = 814.0795000000001
$ws.Range("K131").Value = 1180.5
$ws.Range("L131").Value = 2442.2385
$ws.Range("M131").Value = 3859.5
$ws.Range("N131").Value = -12522.2385
$ws.Range("H135").Value = 2572.4211
$ws.Range("I135").Value = 970.5454999999999
$ws.Range("K135").Value = 8734.9095
$ws.Range("M135").Value = -6199.9095

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 811.8333
$ws.Range("I97").Value = 630
$ws.Range("J97").Value = 902.75
$ws.Range("K97").Value = 630
$ws.Range("L97").Value = 902.75
$ws.Range("M97").Value = -134
$ws.Range("N97").Value = -1894.75

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3046.6365
$ws.Range("I7").Value = 1334.6666
$ws.Range("J7").Value = 5101
$ws.Range("K7").Value = 1334.6666
$ws.Range("L7").Value = 5101
$ws.Range("M7").Value = -1222.6666
$ws.Range("N7").Value = -5325
$ws.Range("H68").Value = 894.5465
$ws.Range("I68").Value = 886.8148
$ws.Range("J68").Value = 1019.8
$ws.Range("K68").Value = 886.8148
$ws.Range("L68").Value = 1019.8
$ws.Range("M68").Value = -137.8148
$ws.Range("N68").Value = -2517.8
$ws.Range("H71").Value = 894.5465
$ws.Range("I71").Value = 886.8148
$ws.Range("J71").Value = 1019.8
$ws.Range("K71").Value = 4434.074
$ws.Range("L71").Value = 5099
$ws.Range("M71").Value = -690.0739999999996
$ws.Range("N71").Value = -12587
$ws.Range("H126").Value = 3046.6365
$ws.Range("I126").Value = 1334.6666
$ws.Range("J126").Value = 5101
$ws.Range("K126").Value = 4003.9998
$ws.Range("L126").Value = 15303
$ws.Range("M126").Value = -1533.9998
$ws.Range("N126").Value = -20243
$ws.Range("H136").Value = 3962.4243
$ws.Range("I136").Value = 1824.4
$ws.Range("J136").Value = 5744.1113
$ws.Range("K136").Value = 5473.200000000001
$ws.Range("L136").Value = 17232.3339
$ws.Range("M136").Value = -2923.200000000001
$ws.Range("N136").Value = -22332.3339

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2719.1155
$ws.Range("I132").Value = 933.25
$ws.Range("J132").Value = 4249.857
$ws.Range("K132").Value = 2799.75
$ws.Range("L132").Value = 12749.571
$ws.Range("M132").Value = -269.75
$ws.Range("N132").Value = -17809.571
